# Auto-generated edit script: refresh market-price derived columns (H-N)
# on the Marilith_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Mirrors a scheduled-runner price refresh: currentAveragePrice* / LevePrice* / LeveProfit*
# values are overwritten in place; a few cells that no longer carry a value are cleared.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = [ordered]@{
    "H19" = 1344.8182
    "J19" = 1932
    "L19" = 1932
    "N19" = -2282
    "H41" = 229.625
    "I41" = 206.16667
    "J41" = 300
    "K41" = 206.16667
    "L41" = 300
    "M41" = 233.83333
    "N41" = -1180
    "H43" = 0
    "I43" = 0
    "J43" = 0
    "K43" = 0
    "L43" = 0
    "M43" = $null
    "N43" = $null
    "H113" = 9927.857
    "I113" = 7373.75
    "J113" = 13333.333
    "K113" = 7373.75
    "L113" = 13333.333
    "M113" = -4119.75
    "N113" = -19841.333
    "H116" = 4121.5713
    "I116" = 3349.25
    "K116" = 3349.25
    "M116" = 92.75
    "H125" = 102811.7
    "I125" = 2041.8572
    "K125" = 18376.7148
    "M125" = -15916.7148
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = [ordered]@{
    "H96" = 18326.334
    "J96" = 18326.334
    "L96" = 18326.334
    "N96" = -23818.334
    "H119" = 200000
    "J119" = 200000
    "L119" = 200000
    "N119" = -209676
    "H122" = 2499.1667
    "I122" = 1999.2
    "J122" = 4999
    "K122" = 5997.6
    "L122" = 14997
    "M122" = -3547.6
    "N122" = -19897
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = [ordered]@{
    "H5" = 2034
    "I5" = 1066.6666
    "J5" = 3001.3333
    "K5" = 1066.6666
    "L5" = 3001.3333
    "M5" = -953.6666
    "N5" = -3227.3333
    "H22" = 125
    "I22" = 100
    "K22" = 100
    "M22" = 73
    "H86" = 699.4
    "I86" = 732.6667
    "J86" = 649.5
    "K86" = 732.6667
    "L86" = 649.5
    "M86" = 390.3333
    "N86" = -2895.5
    "H89" = 699.4
    "I89" = 732.6667
    "J89" = 649.5
    "K89" = 3663.3335
    "L89" = 3247.5
    "M89" = 1952.6665
    "N89" = -14479.5
    "H123" = 0
    "J123" = 0
    "L123" = 0
    "N123" = $null
    "H134" = 5065.273
    "I134" = 4970.3687
    "K134" = 14911.1061
    "M134" = -12376.1061
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("CRP")
$updates = [ordered]@{
    "H2" = 1249.6
    "I2" = 83
    "J2" = 2999.5
    "K2" = 83
    "L2" = 2999.5
    "M2" = 30
    "N2" = -3225.5
    "H16" = 0
    "I16" = 0
    "K16" = 0
    "M16" = $null
    "H58" = 3171.1875
    "I58" = 2671
    "J58" = 3814.2856
    "K58" = 2671
    "L58" = 3814.2856
    "M58" = -2468
    "N58" = -4220.2856
    "H107" = 577.1429000000001
    "I107" = 286.25
    "K107" = 286.25
    "M107" = 1633.75
    "H113" = 0
    "I113" = 0
    "K113" = 0
    "M113" = $null
    "H122" = 1673
    "J122" = 1250
    "L122" = 3750
    "N122" = -8650
    "H132" = 2332
    "I132" = 2332
    "K132" = 6996
    "M132" = -4466
    "H136" = 3171.1875
    "I136" = 2671
    "J136" = 3814.2856
    "K136" = 8013
    "L136" = 11442.8568
    "M136" = -5463
    "N136" = -16542.8568
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = [ordered]@{
    "H34" = 2494.5
    "J34" = 2494.5
    "L34" = 7483.5
    "N34" = -7651.5
    "H39" = 167
    "I39" = 167
    "J39" = 0
    "K39" = 501
    "L39" = 0
    "M39" = -207
    "N39" = $null
    "H46" = 2195
    "I46" = 1325
    "J46" = 3500
    "K46" = 3975
    "L46" = 10500
    "M46" = -3884
    "N46" = -10682
    "H55" = 1217
    "J55" = 2145
    "L55" = 6435
    "N55" = -6789
    "H131" = 0
    "I131" = 0
    "J131" = 0
    "K131" = 0
    "L131" = 0
    "M131" = $null
    "N131" = $null
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = [ordered]@{
    "H7" = 168333.33
    "J7" = 168333.33
    "L7" = 168333.33
    "N7" = -168557.33
    "H8" = 168333.33
    "J8" = 168333.33
    "L8" = 168333.33
    "N8" = -168611.33
    "H10" = 0
    "J10" = 0
    "L10" = 0
    "N10" = $null
    "H102" = 8173.5
    "I102" = 1465.3334
    "K102" = 1465.3334
    "M102" = 156.6666
    "H113" = 3721.75
    "I113" = 2165.6667
    "K113" = 2165.6667
    "M113" = 4.333299999999781
    "H126" = 12256.143
    "I126" = 12256.143
    "K126" = 36768.429
    "M126" = -34298.429
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = [ordered]@{
    "H22" = 1023.4375
    "I22" = 816.4545000000001
    "J22" = 1478.8
    "K22" = 816.4545000000001
    "L22" = 1478.8
    "M22" = -521.4545000000001
    "N22" = -2068.8
    "H27" = 1023.4375
    "I27" = 816.4545000000001
    "J27" = 1478.8
    "K27" = 816.4545000000001
    "L27" = 1478.8
    "M27" = -709.4545000000001
    "N27" = -1692.8
    "H46" = 4374.5
    "I46" = 2998.5
    "J46" = 4833.1665
    "K46" = 2998.5
    "L46" = 4833.1665
    "M46" = -2810.5
    "N46" = -5209.1665
    "H55" = 201.6842
    "I55" = 187
    "J55" = 208.46153
    "K55" = 187
    "L55" = 208.46153
    "M55" = -14
    "N55" = -554.46153
    "H61" = 4130.8335
    "I61" = 6500
    "J61" = 2946.25
    "K61" = 6500
    "L61" = 2946.25
    "M61" = -6298
    "N61" = -3350.25
    "H113" = 4130.8335
    "I113" = 6500
    "J113" = 2946.25
    "K113" = 6500
    "L113" = 2946.25
    "M113" = -4330
    "N113" = -7286.25
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = [ordered]@{
    "H113" = 245.33333
    "I113" = 0
    "J113" = 245.33333
    "K113" = 0
    "L113" = 735.99999
    "M113" = $null
    "N113" = -5075.99999
    "H124" = 0
    "J124" = 0
    "L124" = 0
    "N124" = $null
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
